$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("carte")

# New grid of image references, 5 columns (A..E) x 6 rows
$img1 = "assets/for_tests/img_1.png"
$img2 = "assets/for_tests/img_2.png"
$pika = "assets/for_tests/pikachu.jpeg"

$values = @(
    @($pika, $img1, $img1, $img2, $img1),
    @($img1, $img2, $img1, $img1, $img2),
    @($img1, $img1, $img2, $img1, $img1),
    @($img1, $img2, $img1, $img2, $img1),
    @($img1, $img1, $img2, $img1, $img2)
)

for ($r = 0; $r -lt $values.Length; $r++) {
    $row = $values[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# Row 6 only has D6 and E6 populated
$ws.Cells.Item(6, 4).Value = $img1
$ws.Cells.Item(6, 5).Value = $img1

# Column E width (46 character-units; nudged slightly so the stored
# OOXML width serializes as exactly "46" rather than "46.833...")
$ws.Range("E1:E6").ColumnWidth = 45.14

# Selection
$ws.Range("B9").Select()

$wb.Save()
